# DPLKKLM055-001 - "Update 15 mar 2023"
# - Clears the "run" flag on row 3 (A3)
# - Populates row 4 as a new, active test-data row (mirrors rows 2/3)
# - Moves the sheet view / selection to around AB3 (scrolled to column Q)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 3: the row is no longer flagged to run -> remove the A3 "run" flag
# ---------------------------------------------------------------------
$ws.Range("A3").Clear()

# ---------------------------------------------------------------------
# Row 4: fill in a brand-new test case (same shape as rows 2 & 3)
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 90

$ws.Range("A4").Value = "run"
$ws.Range("B4").Value = "DPLKKLM055-001"
$ws.Range("C4").Value = "Klaim - Transaksi - Calculate Klaim - Send to Approval - Klaim Pasca Kerja"
$ws.Range("D4").Value = "Klaim Pasca Kerja"
$ws.Range("E4").Value = "Calculate Klaim sub send to approval bisa dilakukan dengan baik. Dalam perhitungan nominal klaim, dikenakan biaya administrasi dan fee yang disesuaikan dengan ketentuan PKS"
$ws.Range("G4").Value = 44385
$ws.Range("H4").Value = "bni1234"
$ws.Range("I4").Value = "Klaim"
# J4 already carries the "number stored as text" style (quotePrefix) -
# keep it by writing the value with a leading apostrophe.
$ws.Range("J4").Value = "'Transaksi"
$ws.Range("K4").Value = "Calculate Klaim Pasca Kerja "

# ID Peserta for this row keeps its leading zeroes -> store as text.
$ws.Range("N4").Value = "'0000000045"

# Username/Password/ID Peserta summary formula (same pattern as F2/F3).
$ws.Range("F4").Formula = '= "Username : "&G4&",' + "`n" + 'Password : bni1234,' + "`n" + 'ID Peserta :  "&N4'

# ---------------------------------------------------------------------
# Sheet view: scroll right to column Q and select AB3
# ---------------------------------------------------------------------
$ws.Range("AB3").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 17
$win.ScrollRow = 1
